# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to the freshly-scraped values.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F.
$updates = @{
    2  = 8366
    3  = 7881
    9  = 126
    11 = 231
    12 = 713
    14 = 1872
    17 = 14
    19 = 129
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
